# Add a new "2022" column (L) to the worksheet, mirroring the existing
# 2021 column (K), and move the active selection to L2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# Row 2: empty cell but keep the thick-bottom border style used by the
# rest of that row (same formatting as K2).
$ws.Range("K2").Copy()
$ws.Range("L2").PasteSpecial($xlPasteFormats)

# Row 3: year header, same style as K3.
$ws.Range("K3").Copy()
$ws.Range("L3").PasteSpecial($xlPasteFormats)
$ws.Range("L3").Value = 2022

# Row 4: same style as K4.
$ws.Range("K4").Copy()
$ws.Range("L4").PasteSpecial($xlPasteFormats)
$ws.Range("L4").Value = 370

# Row 5: same style as K5.
$ws.Range("K5").Copy()
$ws.Range("L5").PasteSpecial($xlPasteFormats)
$ws.Range("L5").Value = 137

# Row 6: same style as K6.
$ws.Range("K6").Copy()
$ws.Range("L6").PasteSpecial($xlPasteFormats)
$ws.Range("L6").Value = 314

# Row 7: same style as K7.
$ws.Range("K7").Copy()
$ws.Range("L7").PasteSpecial($xlPasteFormats)
$ws.Range("L7").Value = 121

# Row 8: base formatting copied from K8, then apply its own number format
# (thousands separator) on top, producing a distinct style entry.
$ws.Range("K8").Copy()
$ws.Range("L8").PasteSpecial($xlPasteFormats)
$ws.Range("L8").Value = 50
$ws.Range("L8").NumberFormat = "#,##0"

# Row 9: same style as K9 (bottom row border).
$ws.Range("K9").Copy()
$ws.Range("L9").PasteSpecial($xlPasteFormats)
$ws.Range("L9").Value = 16

$excel.CutCopyMode = $false

# Move the active selection to L2, matching the authored workbook.
$ws.Range("L2").Select()
